$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 413, shifting all existing data (rows 413-443) down by one.
$ws.Rows.Item(413).Insert()

# Populate the newly inserted row 413 with the latest weekly price entry.
$ws.Cells.Item(413, 1).Value = 10
$ws.Cells.Item(413, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(413, 3).Value = "La Araucanía"
$ws.Cells.Item(413, 4).Value = 45013
$ws.Cells.Item(413, 5).Value = 9
$ws.Cells.Item(413, 6).Value = 100112017
$ws.Cells.Item(413, 7).Value = "Apio"
$ws.Cells.Item(413, 8).Value = "Americana (o)"
$ws.Cells.Item(413, 9).Value = "Primera"
$ws.Cells.Item(413, 10).Value = 40
$ws.Cells.Item(413, 11).Value = 10000
$ws.Cells.Item(413, 12).Value = 10000
$ws.Cells.Item(413, 13).Value = 10000
$ws.Cells.Item(413, 14).Value = "$/docena de matas"
$ws.Cells.Item(413, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(413, 16).Value = 1667
$ws.Cells.Item(413, 17).Value = 6
$ws.Cells.Item(413, 18).Value = "Hortaliza"

# Apply the same date-time number format used by the other cells in column D.
$ws.Cells.Item(413, 4).NumberFormat = $ws.Cells.Item(414, 4).NumberFormat
